# Applies the "correct response" column edit for module 2b invitation-task
# mapping sheet: adds an "I" column ("correct") with the letter of the
# answer choice (a/b/c/d) that is correct for each row, and fixes a few
# rows where the choice columns (E:H) were in the wrong order.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------
# 1. Fix the four rows whose E:H (a/b/c/d) choice values were scrambled.
# ---------------------------------------------------------------------
$ws.Cells.Item(4, 5).Value = "10-20"
$ws.Cells.Item(4, 6).Value = "5-10"
$ws.Cells.Item(4, 7).Value = "0-20"
$ws.Cells.Item(4, 8).Value = "0-5"

$ws.Cells.Item(12, 5).Value = "5-5"
$ws.Cells.Item(12, 6).Value = "10-5"
$ws.Cells.Item(12, 7).Value = "10-0"
$ws.Cells.Item(12, 8).Value = "5-0"

$ws.Cells.Item(15, 5).Value = "10-0"
$ws.Cells.Item(15, 6).Value = "0-5"
$ws.Cells.Item(15, 7).Value = "10-5"
$ws.Cells.Item(15, 8).Value = "0-20"

$ws.Cells.Item(27, 5).Value = "10-5"
$ws.Cells.Item(27, 6).Value = "0-2"
$ws.Cells.Item(27, 7).Value = "5-10"
$ws.Cells.Item(27, 8).Value = "10-10"

# ---------------------------------------------------------------------
# 2. Add the new "correct" column (I) with header + per-row answer.
# ---------------------------------------------------------------------
$ws.Cells.Item(1, 9).Value = "correct"

$correct = "d","c","c","a","a","c","b","a","b","a","c","b","d","b","c","a","b","c","a","c","a","a","d","a","c","d","b","a"
for ($i = 0; $i -lt $correct.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $correct[$i]
}

# Match the text formatting used by the other answer columns (E:H).
$ws.Range("I1:I29").NumberFormat = "@"

# ---------------------------------------------------------------------
# 3. Cosmetic touch-ups captured in the diff.
# ---------------------------------------------------------------------
$ws.Columns.Item(4).ColumnWidth = 31.5
$ws.Range("H28").Select() | Out-Null
